# Southwest CMA Compliance Template - "adding functionality to scene session"
#
# The KPIs sheet's "session_level" column (G) carried a "Y" flag on the
# compliance-score KPI rows (18-26). Those flags are being retired, the
# rows get a touch shorter now that the extra column is gone, and the
# KPIs sheet (not Targets) becomes the sheet the workbook opens on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")

# Drop the "session_level" (column G) marker on rows 18-26 and tighten
# the row height now that the cell no longer forces the old wrap height.
for ($row = 18; $row -le 26; $row++) {
    $ws.Range("G$row").ClearContents()
    $ws.Rows.Item($row).RowHeight = 55.2
}

# Make KPIs the active sheet/tab (was "Targets") and leave the cursor on
# the last data row where the removed column used to live.
$ws.Activate()
$ws.Range("G26").Select()
